$d = $word.ActiveDocument

# Helper: split the run covering character offset $endPos.."$pos" (document
# character positions, same numbering as Range.Start/.End) into two runs
# without Word's "identical formatting" run-coalescing silently merging
# them back together. Toggling a character property across the tail
# sub-range (and then reverting it) forces a hard run boundary at $pos
# while leaving the final formatting identical on both sides.
# NB: only plain numbers are passed in (not live Range/COM objects) -
# passing COM objects through function parameters in this interpreter
# does not preserve the live binding needed for the split to stick.
function Split-RunAt($pos, $endPos) {
    $tail = $d.Range($pos, $endPos)
    $tail.Bold = 1
    $tail2 = $d.Range($pos, $endPos)
    $tail2.Bold = 0
}

# Locate the two list-item paragraphs by their current text. Paragraph
# Range.Text includes the trailing paragraph-mark (carriage return), so
# trim it before comparing.
$mp4Para = $null
$readmePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd("`r")
    if ($t -eq "Portfolio_reflection_v4.mp4") { $mp4Para = $para }
    if ($t -eq "README") { $readmePara = $para }
}

# --- "Portfolio_reflection_v4.mp4" -> runs "p" / "ortfolio_reflection_v4.mp4" ---
if ($mp4Para -ne $null) {
    $full = $mp4Para.Range
    $firstChar = $d.Range($full.Start, $full.Start + 1)
    $firstChar.Text = "p"

    $full2 = $mp4Para.Range
    Split-RunAt ($full2.Start + 1) $full2.End
}

# --- "README" -> runs "README" / ".md" ---
if ($readmePara -ne $null) {
    $full = $readmePara.Range
    $end = $d.Range($full.End - 1, $full.End - 1)
    $end.InsertAfter(".md")

    $full2 = $readmePara.Range
    Split-RunAt ($full2.Start + 6) $full2.End
}

Write-Output "mp4=[$($mp4Para.Range.Text)] readme=[$($readmePara.Range.Text)]"
